# Added card number / cost to the card: "条件" (Condition, E2) used to hold
# the cost text "4风"; the cost has grown to "4风1火1水1暗1光" and is now
# also duplicated into "代价" (Cost, H2), which was previously empty.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "4风1火1水1暗1光"
$ws.Range("H2").Value = "4风1火1水1暗1光"

# Move the active selection from J2 to E2 to match the saved view state.
$ws.Range("E2").Select()
